# Update column F ("dSF") values for rows 3-9 and 11 as part of the
# "repull data, push all data, mean calculation" update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = -2
$ws.Range("F4").Value = 4
$ws.Range("F5").Value = -1
$ws.Range("F6").Value = -10
$ws.Range("F7").Value = -1
$ws.Range("F8").Value = 2
$ws.Range("F9").Value = -1
$ws.Range("F11").Value = -7
